$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (A, B, C, E, F, G, H, I - D stays empty)
$ws.Range("A2").Value = 53
$ws.Range("B2").Value = 173
$ws.Range("C2").Value = 152
$ws.Range("E2").Value = 29
$ws.Range("F2").Value = 32
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 33
$ws.Range("I2").Value = 4

# Column C values for rows 3 through 54
$cValues = @(115,84,64,137,36,71,106,83,29,170,106,168,68,82,135,149,118,105,30,31,64,20,7,165,140,67,164,25,99,25,69,131,144,68,121,110,85,72,127,80,152,71,12,168,124,79,26,168,2,76,161,112)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
